# Update crypto price/volume data per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.873.65"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "1.825.47"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'310.91"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "'0.4577"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "'0.3675"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "'0.07156"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'0.8718"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").Value = "'0.07787"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'19.52"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "1.811.37"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "'5.316"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "'6.379"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "'86.81"
$ws.Range("E16").Value = "  -5.51%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.000008693"
$ws.Range("E18").Value = "  -4.14%  "
$ws.Range("D20").Value = "26.891.18"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "'14.44"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "'4.991"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'10.45"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.008"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'150.89"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.17"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'1.952"
$ws.Range("E27").Value = "  -5.83%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "'113.39"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'4.906"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.08804"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("B31").Value = "HuobiToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D31").Value = "'2.995"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.7470"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.468"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.129"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'2.521"
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.084"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01934"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.909"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05105"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.915"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.4961"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1593"
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.261"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.4671"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.007"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.10"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'101.21"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.604"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06087"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'64.32"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.63"
$ws.Range("E51").Value = "  -0.46%  "
